$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Minor timestamp precision correction on existing row 5
$ws.Cells.Item(5, 1).Value = 45878.16685479166

# New row of sensor data (row 6)
$ws.Cells.Item(6, 1).Value = 45878.2085045788
$ws.Cells.Item(6, 1).NumberFormat = $ws.Cells.Item(5, 1).NumberFormat

$ws.Cells.Item(6, 2).Value = 2025
$ws.Cells.Item(6, 3).Value = 37
$ws.Cells.Item(6, 4).Value = 13.17
$ws.Cells.Item(6, 5).Value = 92.56
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 2.37
$ws.Cells.Item(6, 8).Value = "N"
$ws.Cells.Item(6, 9).Value = 0
$ws.Cells.Item(6, 10).Value = "05:00:14"
